$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (not auto-converted to a number/date by
# Excel), while leaving the cell's style index untouched -- ClearFormats
# drops the temporary "@" text number-format right after the value has
# been committed, so no stray style sticks around on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Rows 2-36: refreshed Price (D) / Volume(1h) (E) figures
Set-TextValue $ws.Range("D2") "68.031.61"
$ws.Range("E2").Value = "  -1.52%  "

Set-TextValue $ws.Range("D3") "3.869.80"
$ws.Range("E3").Value = "  -1.56%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue $ws.Range("D5") "599.15"

Set-TextValue $ws.Range("D6") "170.77"
$ws.Range("E6").Value = "  +1.72%  "

Set-TextValue $ws.Range("D7") "3.867.91"
$ws.Range("E7").Value = "  -1.52%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -1.03%  "

$ws.Range("E10").Value = "  -5.14%  "

Set-TextValue $ws.Range("D11") "6.43"
$ws.Range("E11").Value = "  -0.54%  "

Set-TextValue $ws.Range("D12") "0.456"
$ws.Range("E12").Value = "  -2.07%  "

Set-TextValue $ws.Range("D13") "0.0000260"
$ws.Range("E13").Value = "  -0.17%  "

Set-TextValue $ws.Range("D14") "36.97"
$ws.Range("E14").Value = "  -1.77%  "

Set-TextValue $ws.Range("D15") "4.527.01"
$ws.Range("E15").Value = "  -1.37%  "

Set-TextValue $ws.Range("D16") "3.880.91"
$ws.Range("E16").Value = "  -3.21%  "

Set-TextValue $ws.Range("D17") "68.200.23"
$ws.Range("E17").Value = "  -1.35%  "

Set-TextValue $ws.Range("D18") "18.14"
$ws.Range("E18").Value = "  +4.23%  "

Set-TextValue $ws.Range("D19") "7.34"
$ws.Range("E19").Value = "  -2.11%  "

$ws.Range("E20").Value = "  -0.30%  "

Set-TextValue $ws.Range("D21") "10.75"
$ws.Range("E21").Value = "  -2.28%  "

Set-TextValue $ws.Range("D22") "466.92"
$ws.Range("E22").Value = "  -6.13%  "

Set-TextValue $ws.Range("D23") "0.738"
$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("E24").Value = "  -5.62%  "

Set-TextValue $ws.Range("D25") "83.17"

Set-TextValue $ws.Range("D26") "2.23"
$ws.Range("E26").Value = "  -2.18%  "

Set-TextValue $ws.Range("D27") "12.05"
$ws.Range("E27").Value = "  -1.22%  "

$ws.Range("E28").Value = "  +0.01%  "

Set-TextValue $ws.Range("D29") "9.98"
$ws.Range("E29").Value = "  -2.85%  "

Set-TextValue $ws.Range("D30") "2.95"
$ws.Range("E30").Value = "  -1.19%  "

Set-TextValue $ws.Range("D31") "4.020.39"
$ws.Range("E31").Value = "  -1.58%  "

Set-TextValue $ws.Range("D32") "7.73"
$ws.Range("E32").Value = "  -0.80%  "

Set-TextValue $ws.Range("D33") "2.31"
$ws.Range("E33").Value = "  -3.38%  "

Set-TextValue $ws.Range("D34") "31.17"
$ws.Range("E34").Value = "  -2.66%  "

Set-TextValue $ws.Range("D35") "9.49"
$ws.Range("E35").Value = "  -0.14%  "

Set-TextValue $ws.Range("D36") "3.834.50"
$ws.Range("E36").Value = "  -1.63%  "

# Rows 37-38: dogwifhat overtakes Hedera in the ranking, so the two rows
# swap coin/link/price/volume, and Hedera's own figures also refresh.
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D37") "3.75"
$ws.Range("E37").Value = "  +12.94%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D38") "0.104"
$ws.Range("E38").Value = "  -2.72%  "

# Rows 39-51: refreshed Price (D) / Volume(1h) (E) figures
Set-TextValue $ws.Range("D39") "1.02"
$ws.Range("E39").Value = "  -1.95%  "

$ws.Range("E40").Value = "  +0.77%  "

Set-TextValue $ws.Range("D41") "5.90"
$ws.Range("E41").Value = "  -1.76%  "

$ws.Range("E42").Value = "  +0.13%  "

Set-TextValue $ws.Range("D43") "0.312"
$ws.Range("E43").Value = "  -3.15%  "

Set-TextValue $ws.Range("D44") "0.000301"
$ws.Range("E44").Value = "  +4.86%  "

Set-TextValue $ws.Range("D45") "424.00"
$ws.Range("E45").Value = "  -2.22%  "

Set-TextValue $ws.Range("D46") "1.97"
$ws.Range("E46").Value = "  -1.80%  "

$ws.Range("E47").Value = "  +0.00%  "

Set-TextValue $ws.Range("D48") "8.62"

Set-TextValue $ws.Range("D49") "47.11"
$ws.Range("E49").Value = "  -1.83%  "

Set-TextValue $ws.Range("D50") "26.76"
$ws.Range("E50").Value = "  +4.01%  "

Set-TextValue $ws.Range("D51") "142.92"
$ws.Range("E51").Value = "  -0.11%  "
